$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Close out the "2018.11.28 第十三周周三" week: mark 陈升云 (row 217) and
#    吴帅辰 (row 219) completion fraction as 100% (1).
# ---------------------------------------------------------------------------
$ws.Range("C217").Value = 1
$ws.Range("C219").Value = 1

# ---------------------------------------------------------------------------
# 2) Append a new weekly block for "2018.12.5 第十四周周三" in rows 224-232,
#    mirroring the structure of the previous block (rows 214-222):
#      224      -> merged date header (A224:D224)
#      225      -> column headers (组员/计划内容/完成情况/备注)
#      226-230  -> one row per team member
#      231-232  -> merged "总结：" summary row (A231:D232)
# ---------------------------------------------------------------------------

# -- values -------------------------------------------------------------
$ws.Range("A224").Value = "日期：2018.12.5 第十四周周三"

$ws.Range("A225").Value = "组员"
$ws.Range("B225").Value = "计划内容"
$ws.Range("C225").Value = "完成情况"
$ws.Range("D225").Value = "备注"

$ws.Range("A226").Value = "王伟锋"
$ws.Range("B226").Value = "基本功能已完成，协助app完成"
$ws.Range("D226").Value = "协助情况不计入任务完成情况"

$ws.Range("A227").Value = "陈升云"
$ws.Range("B227").Value = "基本功能已完成，开始优化客户端和解决bug"

$ws.Range("A228").Value = "林玮成"
$ws.Range("B228").Value = "等待最终程序的测试"

$ws.Range("A229").Value = "吴帅辰"
$ws.Range("B229").Value = "基本功能已完成，协助app完成"
$ws.Range("D229").Value = "协助情况不计入任务完成情况"

$ws.Range("A230").Value = "李海洋"
$ws.Range("B230").Value = "基本功能已完成，开始优化客户端和解决bug"

$ws.Range("A231").Value = "总结："

# -- formatting: copy from the matching cell of the previous block ------
function CopyFmt($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

CopyFmt "A214" "A224"
CopyFmt "B214" "B224"
CopyFmt "C214" "C224"
CopyFmt "D214" "D224"

CopyFmt "A215" "A225"
CopyFmt "B215" "B225"
CopyFmt "C215" "C225"
CopyFmt "D215" "D225"

CopyFmt "A216" "A226"
CopyFmt "B216" "B226"
CopyFmt "C216" "C226"
CopyFmt "D216" "D226"

CopyFmt "A217" "A227"
CopyFmt "B217" "B227"
CopyFmt "C217" "C227"
CopyFmt "D217" "D227"

CopyFmt "A218" "A228"
CopyFmt "B218" "B228"
CopyFmt "C218" "C228"
CopyFmt "D218" "D228"

CopyFmt "A219" "A229"
CopyFmt "B219" "B229"
CopyFmt "C219" "C229"
CopyFmt "D219" "D229"

CopyFmt "A220" "A230"
CopyFmt "B220" "B230"
CopyFmt "C220" "C230"
CopyFmt "D220" "D230"

CopyFmt "A221" "A231"
CopyFmt "B221" "B231"
CopyFmt "C221" "C231"
CopyFmt "D221" "D231"

CopyFmt "A222" "A232"
CopyFmt "B222" "B232"
CopyFmt "C222" "C232"
CopyFmt "D222" "D232"

# -- merges ---------------------------------------------------------------
$ws.Range("A224:D224").Merge()
$ws.Range("A231:D232").Merge()

# ---------------------------------------------------------------------------
# 3) Scroll the view to show the newly-added rows, like the source author's
#    last saved position.
# ---------------------------------------------------------------------------
$ws.Range("F235").Select()
$excel.ActiveWindow.ScrollRow = 199
